# Update the "想去人数" (want-to-go count) figures that were refreshed when
# the gh-pages data generation job re-ran (commit 456a3b4).
#
# Sheet "展览" (Exhibition) and sheet "全部类型" (All types) both contain the
# same rows of data; update column F (想去人数) for rows 2, 4, 5, 6 and 7 on
# each sheet. Note the two sheets diverge for row 5: "展览" picks up the
# freshly scraped value (7711) while "全部类型" ends up with 0 for that row.

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F2").Value = 2253
$wsExhibition.Range("F4").Value = 1664
$wsExhibition.Range("F5").Value = 7711
$wsExhibition.Range("F6").Value = 181
$wsExhibition.Range("F7").Value = 219

$wsAllTypes = $wb.Worksheets.Item("全部类型")
$wsAllTypes.Range("F2").Value = 2253
$wsAllTypes.Range("F4").Value = 1664
$wsAllTypes.Range("F5").Value = 0
$wsAllTypes.Range("F6").Value = 181
$wsAllTypes.Range("F7").Value = 219
